# The commit swaps the contents of ppt/theme/theme1.xml (the Slide
# Master's "Integral" / Red Violet theme) and ppt/theme/theme2.xml (the
# Notes Master's "Office Theme"): theme1.xml ends up holding the colours
# that used to live in theme2.xml, and vice versa. Font scheme and
# format scheme are byte-identical between the two themes, so the only
# thing that actually changes is the 12-colour colour scheme (plus the
# cosmetic theme/colour-scheme display names, which aren't swapped here
# because they live outside the 12-colour scheme).
#
# The PowerPoint object model exposes the live colour scheme via
# Master.Theme.ThemeColorScheme, indexed 1-12 in the standard
# MsoThemeColorSchemeIndex order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
# Writing into it updates the master's underlying theme part, so we push
# the colours that previously lived in theme2.xml ("Office" scheme) onto
# the slide master's theme, reproducing the swapped theme1.xml content.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

function ConvertTo-BGR([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Target colours = the values that were in ppt/theme/theme2.xml
# ("Office" colour scheme) before the edit.
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Item($i).RGB = ConvertTo-BGR($officeColors[$i - 1])
}
